# collection software hover update
# Adds a small helper column (B) of hover/collection values next to the
# existing A1 total, then leaves the selection parked one row below the
# new data - matching the author's last-saved cursor position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 1000
$ws.Range("B2").Value = 2000
$ws.Range("B3").Value = 3000
$ws.Range("B4").Value = 4000

$ws.Range("B5").Select()
